$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SinhVien")

# E3: set password value for the newly added student
$ws.Range("E3").Value = "wsq"

# F3: registration date - drop the time-of-day fraction, keep same date (2020-10-16)
$ws.Range("F3").Value = 44120

# H3: gender changed from NAM to NU
$ws.Range("H3").Value = "NU"
